# Weekly driver report update for 2025-04-21
# Update Critical Minutes / Good Roaming Calculation / Totals figures
# on the "Driver Summary" sheet to reflect the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2
$ws.Range("C3").Value = 2868
$ws.Range("D3").Value = 88.3

# Row 4: MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1297
$ws.Range("C4").Value = 206

# Row 5: Totals
$ws.Range("C5").Value = 3074
